$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows for "git rm" (typed A, C, B)
$ws.Range("A9").Value = "git rm"
$ws.Range("C9").Value = "git rm file1.txt"
$ws.Range("B9").Value = "remove a file from the filesystem"

# New rows for "git rm --cached" (typed A, C, B)
$ws.Range("A10").Value = "git rm --cached"
$ws.Range("C10").Value = "git rm --cached  file1.txt"
$ws.Range("B10").Value = "remove a file from the git repository, but not from the filesystem. File becomes untracked"

# Add the missing Example for the "echo" row (row 5)
$ws.Range("C5").Value = 'echo "Hello, world"'

# New rows for "echo ... > file" (write content) (typed A, B, C)
$ws.Range("A11").Value = 'echo "content to write to file" > file_name.txt'
$ws.Range("B11").Value = "write the content to the file, create the file at the mean time"
$ws.Range("C11").Value = 'echo "hello world" > another-file.txt'

# New rows for "ls" (typed A, B)
$ws.Range("A12").Value = "ls"
$ws.Range("B12").Value = "list files in a folder"

# New rows for "ls -la" (typed A, B)
$ws.Range("A13").Value = "ls -la"
$ws.Range("B13").Value = "list files in a folder, including hiden ones"

# New rows for "cat" (typed A, B, C)
$ws.Range("A14").Value = "cat file_name"
$ws.Range("B14").Value = "print the content of the file in the termimal"
$ws.Range("C14").Value = "cat another-file.txt"

# Column B widened to fit new, longer description text (auto best-fit for the
# new, longer "remove a file from the git repository..." description)
$ws.Columns.Item(2).ColumnWidth = 78.25

# Update selection on the sheet to match post-edit state
$ws.Range("C18").Select()
